$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the date as a formula returning a text literal, then freeze it to a
# plain value via copy/paste-special. This avoids Excel's automatic
# text->date conversion (which would store a date serial + number format)
# while also avoiding any extra cell style (e.g. quote-prefix / text format)
# being allocated, so the cell ends up as a plain shared-string cell just
# like the existing rows.
$ws.Range("A11").Formula = '="1/24/2010"'
$ws.Range("A11").Copy()
$ws.Range("A11").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B11").Value = 1.75
$ws.Range("C11").Value = "CFP Update/Edit"

$ws.Range("A12").Select()
